# Updates the PSSM (Position-Specific Scoring Matrix) values on Sheet1
# with refreshed results from the supplemental-figure re-run.
# Rows 2-21 correspond to the 20 amino-acid residues (column A), and
# columns B-K hold the scores for sequence positions 0-9 (row 1 headers).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New PSSM values, row-major, one inner array per worksheet row (2..21),
# each holding the 10 column values for B..K in order.
$pssmValues = @(
    @(-18.92812797035466, 2.404295995825178, -18.92812797035466, -18.92812797035466, -18.92812797035466, -18.92812797035466, -18.92812797035466, -18.92812797035466, -18.92812797035466, -18.92812797035466),
    @(-18.92812797035466, -18.92812797035466, -18.92812797035466, -18.92812797035466, -18.92812797035466, -18.92812797035466, -18.92812797035466, 2.3224347332845, -18.92812797035466, -18.92812797035466),
    @(-18.92812797035466, 2.19193077925869, 2.81343982328956, -18.92812797035466, 2.497436781149565, -18.92812797035466, 4.321925347219229, -18.92812797035466, 2.629506338900482, -18.92812797035466),
    @(-18.92812797035466, 0.941148067627292, -18.92812797035466, -18.92812797035466, -18.92812797035466, 1.962781219587099, -18.92812797035466, -18.92812797035466, -18.92812797035466, -18.92812797035466),
    @(-18.92812797035466, -18.92812797035466, -18.92812797035466, -18.92812797035466, -18.92812797035466, -18.92812797035466, -18.92812797035466, -18.92812797035466, -18.92812797035466, -18.92812797035466),
    @(2.966767116593306, -18.92812797035466, -18.92812797035466, -18.92812797035466, -18.92812797035466, -18.92812797035466, -18.92812797035466, -18.92812797035466, -18.92812797035466, -18.92812797035466),
    @(-18.92812797035466, -18.92812797035466, -18.92812797035466, 2.836339117416804, -18.92812797035466, -18.92812797035466, -18.92812797035466, -18.92812797035466, -18.92812797035466, -18.92812797035466),
    @(3.60669566738692, -18.92812797035466, -18.92812797035466, -18.92812797035466, -18.92812797035466, -18.92812797035466, -18.92812797035466, -18.92812797035466, -18.92812797035466, -18.92812797035466),
    @(-18.92812797035466, -18.92812797035466, -18.92812797035466, -18.92812797035466, -18.92812797035466, -18.92812797035466, -18.92812797035466, 1.607171140603865, -18.92812797035466, 2.278984682135621),
    @(-18.92812797035466, -18.92812797035466, -18.92812797035466, 1.976663185251791, -18.92812797035466, 2.544876635690116, -18.92812797035466, -18.92812797035466, -18.92812797035466, 1.409782429604073),
    @(-18.92812797035466, -18.92812797035466, -18.92812797035466, -18.92812797035466, -18.92812797035466, -18.92812797035466, -18.92812797035466, -18.92812797035466, -18.92812797035466, -18.92812797035466),
    @(-18.92812797035466, -18.92812797035466, -18.92812797035466, 1.726029903022233, -18.92812797035466, -18.92812797035466, -18.92812797035466, -18.92812797035466, 2.248011691543461, 1.697795755777012),
    @(-18.92812797035466, -18.92812797035466, 1.656740051470088, -18.92812797035466, -18.92812797035466, -18.92812797035466, -18.92812797035466, -18.92812797035466, -18.92812797035466, 2.085830502375415),
    @(-18.92812797035466, -18.92812797035466, -0.3592329459695739, -18.92812797035466, -18.92812797035466, -18.92812797035466, -18.92812797035466, -18.92812797035466, -18.92812797035466, -18.92812797035466),
    @(-18.92812797035466, -18.92812797035466, -18.92812797035466, -18.92812797035466, -18.92812797035466, -18.92812797035466, -18.92812797035466, -18.92812797035466, 2.290771631902646, -18.92812797035466),
    @(-18.92812797035466, 0.6079313821857513, -0.08196120228935534, -18.92812797035466, -18.92812797035466, -18.92812797035466, -18.92812797035466, 1.016911434522972, 1.136426668117071, -18.92812797035466),
    @(-18.92812797035466, -18.92812797035466, -18.92812797035466, -18.92812797035466, -18.92812797035466, -18.92812797035466, -18.92812797035466, 0.7352511895555389, 0.9782063031812057, -18.92812797035466),
    @(-18.92812797035466, -18.92812797035466, 1.764505408280157, -18.92812797035466, -18.92812797035466, -18.92812797035466, -18.92812797035466, 2.137017933484948, -18.92812797035466, -18.92812797035466),
    @(-18.92812797035466, 1.721195754712368, 2.231256210073635, -18.92812797035466, 3.843295607427653, -18.92812797035466, -18.92812797035466, 1.95042666818743, -18.92812797035466, 2.322041467788659),
    @(-18.92812797035466, 1.763756089437447, -18.92812797035466, 2.489012203627472, -18.92812797035466, 3.359838267372936, -18.92812797035466, -18.92812797035466, -18.92812797035466, -18.92812797035466)
)

$startRow = 2
$startCol = 2  # column B

for ($i = 0; $i -lt $pssmValues.Length; $i++) {
    $rowValues = $pssmValues[$i]
    $row = $startRow + $i
    for ($j = 0; $j -lt $rowValues.Length; $j++) {
        $col = $startCol + $j
        $ws.Cells.Item($row, $col).Value = $rowValues[$j]
    }
}
